# Applies the "Updated symbol list" crypto data refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apostrophe prefix forces Excel to keep numeric-looking strings as text,
# matching the workbook's original inlineStr / shared-string cell layout.
$q = "'"

$ws.Range("D2").Value = $q + '265.90'
$ws.Range("G2").Value = $q + '10'
$ws.Range("D3").Value = $q + '22.68'
$ws.Range("G3").Value = $q + '10'
$ws.Range("D4").Value = $q + '6.197'
$ws.Range("G4").Value = $q + '10'
$ws.Range("G5").Value = $q + '10'
$ws.Range("D6").Value = $q + '3.559'
$ws.Range("G6").Value = $q + '10'
$ws.Range("D7").Value = $q + '6.708'
$ws.Range("G7").Value = $q + '10'
$ws.Range("G8").Value = $q + '10'
$ws.Range("D9").Value = $q + '0.8248'
$ws.Range("G9").Value = $q + '10'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = $q + '0.1596'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("G10").Value = $q + '10'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = $q + '0.08203'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("G11").Value = $q + '10'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = $q + '0.03405'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G12").Value = $q + '10'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = $q + '0.03155'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("G13").Value = $q + '10'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = $q + '0.09237'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("G14").Value = $q + '10'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = $q + '3.900'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("G15").Value = $q + '10'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = $q + '0.001709'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("G16").Value = $q + '10'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = $q + '0.04830'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("G17").Value = $q + '10'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = $q + '0.0006226'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("G18").Value = $q + '10'
$ws.Range("D19").Value = $q + '0.006174'
$ws.Range("G19").Value = $q + '10'
$ws.Range("D20").Value = $q + '0.006271'
$ws.Range("G20").Value = $q + '10'
$ws.Range("D21").Value = $q + '0.001098'
$ws.Range("G21").Value = $q + '10'
$ws.Range("D22").Value = $q + '0.0001499'
$ws.Range("G22").Value = $q + '10'
$ws.Range("D23").Value = $q + '3.704'
$ws.Range("G23").Value = $q + '10'
$ws.Range("D24").Value = $q + '2.263'
$ws.Range("G24").Value = $q + '10'
$ws.Range("D25").Value = $q + '0.3382'
$ws.Range("G25").Value = $q + '10'
$ws.Range("D26").Value = $q + '0.1207'
$ws.Range("G26").Value = $q + '10'
$ws.Range("D27").Value = $q + '0.0002680'
$ws.Range("G27").Value = $q + '10'
$ws.Range("G28").Value = $q + '10'
$ws.Range("G29").Value = $q + '10'
$ws.Range("G30").Value = $q + '10'
$ws.Range("G31").Value = $q + '10'
$ws.Range("G32").Value = $q + '10'
$ws.Range("G33").Value = $q + '10'
$ws.Range("G34").Value = $q + '10'
$ws.Range("G35").Value = $q + '10'
$ws.Range("G36").Value = $q + '10'
$ws.Range("G37").Value = $q + '10'
$ws.Range("G38").Value = $q + '10'
$ws.Range("G39").Value = $q + '10'
$ws.Range("D40").Value = $q + '0.04591'
$ws.Range("G40").Value = $q + '10'
$ws.Range("D41").Value = $q + '0.006987'
$ws.Range("G41").Value = $q + '10'
$ws.Range("G42").Value = $q + '10'
$ws.Range("D43").Value = $q + '0.003129'
$ws.Range("G43").Value = $q + '10'
$ws.Range("D44").Value = $q + '0.01069'
$ws.Range("G44").Value = $q + '10'
$ws.Range("D45").Value = $q + '0.00006113'
$ws.Range("G45").Value = $q + '10'
$ws.Range("G46").Value = $q + '10'
$ws.Range("D47").Value = $q + '0.7696'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("G47").Value = $q + '10'
$ws.Range("D48").Value = $q + '0.2063'
$ws.Range("G48").Value = $q + '10'
$ws.Range("D49").Value = $q + '0.00002099'
$ws.Range("G49").Value = $q + '10'
$ws.Range("D50").Value = $q + '0.01239'
$ws.Range("G50").Value = $q + '10'
$ws.Range("G51").Value = $q + '10'
